$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 10268.6
$ws.Range("I6").Value = 12585.75
$ws.Range("K6").Value = 37757.25
$ws.Range("M6").Value = -37645.25
# Row 39
$ws.Range("H39").Value = 689.75
$ws.Range("I39").Value = 898.55554
$ws.Range("J39").Value = 63.333332
$ws.Range("K39").Value = 2695.66662
$ws.Range("L39").Value = 189.999996
$ws.Range("M39").Value = -2399.66662
$ws.Range("N39").Value = -781.999996
# Row 53
$ws.Range("H53").Value = 864.6667
$ws.Range("I53").Value = 864.6667
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 864.6667
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -227.6667
$ws.Range("N53").ClearContents()
# Row 61
$ws.Range("H61").Value = 130.5
$ws.Range("I61").Value = 130.5
$ws.Range("K61").Value = 391.5
$ws.Range("M61").Value = -219.5
# Row 64
$ws.Range("H64").Value = 3684.7058
$ws.Range("I64").Value = 3648
$ws.Range("J64").Value = 3737.1428
$ws.Range("K64").Value = 3648
$ws.Range("L64").Value = 3737.1428
$ws.Range("M64").Value = -3400
$ws.Range("N64").Value = -4233.1428
# Row 67
$ws.Range("H67").Value = 3684.7058
$ws.Range("I67").Value = 3648
$ws.Range("J67").Value = 3737.1428
$ws.Range("K67").Value = 3648
$ws.Range("L67").Value = 3737.1428
$ws.Range("M67").Value = -2790
$ws.Range("N67").Value = -5453.1428
# Row 104
$ws.Range("H104").Value = 667.5454999999999
$ws.Range("I104").Value = 667.5454999999999
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 2002.6365
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -255.6364999999998
$ws.Range("N104").ClearContents()
# Row 112
$ws.Range("H112").Value = 2638.5557
$ws.Range("J112").Value = 3192.6428
$ws.Range("L112").Value = 9577.928400000001
$ws.Range("N112").Value = -11793.9284
# Row 116
$ws.Range("H116").Value = 3581.8462
$ws.Range("I116").Value = 2911.2
$ws.Range("K116").Value = 2911.2
$ws.Range("M116").Value = 530.8000000000002
# Row 129
$ws.Range("H129").Value = 768.0769
$ws.Range("J129").Value = 854
$ws.Range("L129").Value = 2562
$ws.Range("N129").Value = -12562
# Row 135
$ws.Range("H135").Value = 32258974
$ws.Range("I135").Value = 471.36365
$ws.Range("J135").Value = 111113090
$ws.Range("K135").Value = 4242.27285
$ws.Range("L135").Value = 1000017810
$ws.Range("M135").Value = -1707.27285
$ws.Range("N135").Value = -1000022880
# Row 137
$ws.Range("H137").Value = 1204.238
$ws.Range("I137").Value = 946.7143
$ws.Range("J137").Value = 1719.2858
$ws.Range("K137").Value = 2840.1429
$ws.Range("L137").Value = 5157.857400000001
$ws.Range("M137").Value = -290.1428999999998
$ws.Range("N137").Value = -10257.8574

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7323.384
$ws.Range("I32").Value = 5993.2876
$ws.Range("K32").Value = 5993.2876
$ws.Range("M32").Value = -5706.2876
# Row 74
$ws.Range("H74").Value = 3003.4375
$ws.Range("I74").Value = 2476.4546
$ws.Range("J74").Value = 4162.8
$ws.Range("K74").Value = 2476.4546
$ws.Range("L74").Value = 4162.8
$ws.Range("M74").Value = -1602.4546
$ws.Range("N74").Value = -5910.8
# Row 77
$ws.Range("H77").Value = 3003.4375
$ws.Range("I77").Value = 2476.4546
$ws.Range("J77").Value = 4162.8
$ws.Range("K77").Value = 12382.273
$ws.Range("L77").Value = 20814
$ws.Range("M77").Value = -8014.273000000001
$ws.Range("N77").Value = -29550
# Row 122
$ws.Range("H122").Value = 2450.2666
$ws.Range("I122").Value = 2016.4445
$ws.Range("J122").Value = 3101
$ws.Range("K122").Value = 6049.333500000001
$ws.Range("L122").Value = 9303
$ws.Range("M122").Value = -3599.333500000001
$ws.Range("N122").Value = -14203

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
# Row 94
$ws.Range("H94").Value = 7353212.5
$ws.Range("I94").Value = 8620985
$ws.Range("J94").Value = 134.8
$ws.Range("K94").Value = 8620985
$ws.Range("L94").Value = 134.8
$ws.Range("M94").Value = -8620534
$ws.Range("N94").Value = -1036.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1333.907
$ws.Range("I31").Value = 1333.907
$ws.Range("K31").Value = 1333.907
$ws.Range("M31").Value = -1038.907
# Row 34
$ws.Range("H34").Value = 1333.907
$ws.Range("I34").Value = 1333.907
$ws.Range("K34").Value = 1333.907
$ws.Range("M34").Value = -1131.907
# Row 107
$ws.Range("H107").Value = 741.4583
$ws.Range("I107").Value = 410.26315
$ws.Range("K107").Value = 410.26315
$ws.Range("M107").Value = 1509.73685
# Row 132
$ws.Range("H132").Value = 2168.6206
$ws.Range("I132").Value = 1814.875
$ws.Range("J132").Value = 2604
$ws.Range("K132").Value = 5444.625
$ws.Range("L132").Value = 7812
$ws.Range("M132").Value = -2914.625
$ws.Range("N132").Value = -12872

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 45003330
$ws.Range("I70").Value = 41670144
$ws.Range("J70").Value = 50003104
$ws.Range("K70").Value = 41670144
$ws.Range("L70").Value = 50003104
$ws.Range("M70").Value = -41669874
$ws.Range("N70").Value = -50003644
# Row 73
$ws.Range("H73").Value = 45003330
$ws.Range("I73").Value = 41670144
$ws.Range("J73").Value = 50003104
$ws.Range("K73").Value = 41670144
$ws.Range("L73").Value = 50003104
$ws.Range("M73").Value = -41669208
$ws.Range("N73").Value = -50004976
# Row 102
$ws.Range("H102").Value = 1491.4286
$ws.Range("I102").Value = 1414.2963
$ws.Range("J102").Value = 1751.75
$ws.Range("K102").Value = 1414.2963
$ws.Range("L102").Value = 1751.75
$ws.Range("M102").Value = 207.7037
$ws.Range("N102").Value = -4995.75
# Row 112
$ws.Range("H112").Value = 34643.7
$ws.Range("J112").Value = 34643.7
$ws.Range("L112").Value = 34643.7
$ws.Range("N112").Value = -36859.7

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5568.077
$ws.Range("I40").Value = 2643
$ws.Range("J40").Value = 7396.25
$ws.Range("K40").Value = 2643
$ws.Range("L40").Value = 7396.25
$ws.Range("M40").Value = -2507
$ws.Range("N40").Value = -7668.25
# Row 61
$ws.Range("H61").Value = 976.2
$ws.Range("I61").Value = 846.8333
$ws.Range("K61").Value = 846.8333
$ws.Range("M61").Value = -644.8333
# Row 93
$ws.Range("H93").Value = 1007.3571
$ws.Range("I93").Value = 961.7692
$ws.Range("K93").Value = 961.7692
$ws.Range("M93").Value = 286.2308
# Row 100
$ws.Range("H100").Value = 1053.3334
$ws.Range("I100").Value = 909.0909
$ws.Range("J100").Value = 1450
$ws.Range("K100").Value = 909.0909
$ws.Range("L100").Value = 1450
$ws.Range("M100").Value = -368.0909
$ws.Range("N100").Value = -2532
# Row 113
$ws.Range("H113").Value = 976.2
$ws.Range("I113").Value = 846.8333
$ws.Range("K113").Value = 846.8333
$ws.Range("M113").Value = 1323.1667
# Row 122
$ws.Range("H122").Value = 62501748
$ws.Range("J122").Value = 2502.5
$ws.Range("L122").Value = 7507.5
$ws.Range("N122").Value = -12407.5
# Row 136
$ws.Range("H136").Value = 2442.6667
$ws.Range("I136").Value = 2046
$ws.Range("J136").Value = 2760
$ws.Range("K136").Value = 6138
$ws.Range("L136").Value = 8280
$ws.Range("M136").Value = -3588
$ws.Range("N136").Value = -13380

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 108
$ws.Range("H108").Value = 32906.25
$ws.Range("J108").Value = 32906.25
$ws.Range("L108").Value = 32906.25
$ws.Range("N108").Value = -40586.25
# Row 113
$ws.Range("H113").Value = 512.0769
$ws.Range("I113").Value = 381.5
$ws.Range("J113").Value = 721
$ws.Range("K113").Value = 1144.5
$ws.Range("L113").Value = 2163
$ws.Range("M113").Value = 1025.5
$ws.Range("N113").Value = -6503
# Row 126
$ws.Range("H126").Value = 45455572
$ws.Range("I126").Value = 66667172
$ws.Range("K126").Value = 200001516
$ws.Range("M126").Value = -199999046
# Row 137
$ws.Range("H137").Value = 32822.855
$ws.Range("J137").Value = 32822.855
$ws.Range("L137").Value = 32822.855
$ws.Range("N137").Value = -43022.855
